$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3076540231704712
$ws.Range("B1").Value = 0.2311508804559708
$ws.Range("C1").Value = 0.3323621153831482
$ws.Range("D1").Value = 3.724145174026489
$ws.Range("E1").Value = 1.613147020339966
